$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '33.766.54'
$ws.Range('E2').Value = '  +7.80%  '
$ws.Range('D3').Value = '1.776.98'
$ws.Range('E3').Value = '  +4.18%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '''225.30'
$ws.Range('E5').Value = '  +1.75%  '
$ws.Range('D6').Value = '''0.560'
$ws.Range('E6').Value = '  +4.72%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '''30.49'
$ws.Range('E8').Value = '  +2.04%  '
$ws.Range('D9').Value = '''46.58'
$ws.Range('E9').Value = '  +4.02%  '
$ws.Range('E10').Value = '  +3.57%  '
$ws.Range('D11').Value = '''0.0666'
$ws.Range('E11').Value = '  +3.64%  '
$ws.Range('E12').Value = '  +1.37%  '
$ws.Range('D13').Value = '2.031.34'
$ws.Range('E13').Value = '  +4.22%  '
$ws.Range('D14').Value = '1.777.16'
$ws.Range('E14').Value = '  +4.05%  '
$ws.Range('D15').Value = '''0.626'
$ws.Range('E15').Value = '  +2.40%  '
$ws.Range('D16').Value = '33.746.46'
$ws.Range('E16').Value = '  +7.82%  '
$ws.Range('D17').Value = '''10.04'
$ws.Range('E17').Value = '  -2.02%  '
$ws.Range('D18').Value = '''4.18'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('D19').Value = '''68.54'
$ws.Range('E19').Value = '  +2.22%  '
$ws.Range('D20').Value = '''251.90'
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('D21').Value = '0.0₃0738'
$ws.Range('E21').Value = '  +1.81%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('E24').Value = '  -2.43%  '
$ws.Range('D25').Value = '''2.14'
$ws.Range('E25').Value = '  -1.20%  '
$ws.Range('D26').Value = '''159.40'
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('D27').Value = '''16.49'
$ws.Range('E27').Value = '  +2.85%  '
$ws.Range('E28').Value = '  +1.09%  '
$ws.Range('D29').Value = '''6.94'
$ws.Range('E29').Value = '  +2.00%  '
$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').Value = '''0.0513'
$ws.Range('E32').Value = '  +1.73%  '
$ws.Range('D33').Value = '''1.20'
$ws.Range('E33').Value = '  +3.89%  '
$ws.Range('D34').Value = '''3.55'
$ws.Range('E34').Value = '  +3.98%  '
$ws.Range('E35').Value = '  +7.78%  '
$ws.Range('D36').Value = '1.481.35'
$ws.Range('E36').Value = '  -1.90%  '
$ws.Range('E37').Value = '  +2.95%  '
$ws.Range('D38').Value = '''0.633'
$ws.Range('E38').Value = '  +3.30%  '
$ws.Range('E39').Value = '  +2.41%  '
$ws.Range('D40').Value = '''83.07'
$ws.Range('E40').Value = '  -0.22%  '
$ws.Range('E41').Value = '  +2.14%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').Value = '''0.885'
$ws.Range('E43').Value = '  +3.58%  '
$ws.Range('E44').Value = '  +1.96%  '
$ws.Range('E45').Value = '  +1.59%  '
$ws.Range('E46').Value = '  +3.85%  '
$ws.Range('D47').Value = '1.930.16'
$ws.Range('E47').Value = '  +4.92%  '
$ws.Range('D48').Value = '''5.72'
$ws.Range('E48').Value = '  +2.39%  '
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').Value = '''11.84'
$ws.Range('E50').Value = '  +14.03%  '
$ws.Range('D51').Value = '''50.74'
$ws.Range('E51').Value = '  -2.93%  '
